$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 546; existing rows 546:565 shift down to 547:566.
$ws.Rows("546:546").Insert()

# Populate the newly inserted row with the new data record.
$ws.Range("A546").Value = 10
$ws.Range("B546").Value = "Vega Modelo de Temuco"
$ws.Range("C546").Value = "La Araucanía"
$ws.Range("D546").Value = 45075
$ws.Range("E546").Value = 9
$ws.Range("F546").Value = 100112040
$ws.Range("G546").Value = "Cilantro"
$ws.Range("H546").Value = "Sin especificar"
$ws.Range("I546").Value = "Primera"
$ws.Range("J546").Value = 55
$ws.Range("K546").Value = 4000
$ws.Range("L546").Value = 4000
$ws.Range("M546").Value = 4000
$ws.Range("N546").Value = "$/docena de atados (2 kilos)"
$ws.Range("O546").Value = "Región Metropolitana"
$ws.Range("P546").Value = 2000
$ws.Range("Q546").Value = 2
$ws.Range("R546").Value = "Hortaliza"
